$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the new commit entry on row 36 (previously blank)
$ws.Range("C36").Value = "external source update & bugs fixed"
$ws.Range("G36").Value = 3

# 2. Insert 10 new blank rows after row 37 (rows 38-47), pushing the
#    old row 39 (Total line) down to row 49, leaving row 48 blank/empty
#    just like row 38 used to be blank before row 39 in the original file.
$ws.Rows("39:48").Insert()

# 3. Give the newly inserted rows (38-47) the same per-column formatting
#    that row 37 already had (C: "20% - Enfasis5" style, D-F: "20% -
#    Enfasis5" style, G: "Neutral" style) so they visually match the
#    other blank rows in the table.
for ($r = 38; $r -le 47; $r++) {
    $ws.Range("C$r").Style = "20% - Énfasis5"
    $ws.Range("D$r").Style = "20% - Énfasis5"
    $ws.Range("E$r").Style = "20% - Énfasis5"
    $ws.Range("F$r").Style = "20% - Énfasis5"
    $ws.Range("G$r").Style = "Neutral"
}

# 4. Apply a single-underline font to I41 (an otherwise empty cell that
#    picked up formatting in the source edit).
$ws.Range("I41").Font.Underline = 2

# 5. Update the Total formula so it includes the newly-filled row 36 and
#    recompute its cached value.
$ws.Range("G49").Formula = "=SUM(G4:G36)"

# 6. Restore the view state recorded in the edited workbook.
$ws.Application.ActiveWindow.ScrollRow = 34
$ws.Range("C36").Select()
